$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.540.63"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -4.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.488.52"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -5.69%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "543.67"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.93"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -4.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.578"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.515.27"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -4.60%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.94%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.48"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.965.90"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "24.55"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "59.672.80"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -4.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000140"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.506.88"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -5.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.38"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.37"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "326.54"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.995"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.82"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.57"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.447"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -10.70%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.00%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.80"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0796"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.28"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.78"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.09%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "158.56"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.44"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "19.04"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.49"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.86%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.73"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.94"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "315.44"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -5.29%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.75"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.65%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.79"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.833"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -7.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.995"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.603"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.71"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "126.65"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0534"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0940"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0233"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.71"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -5.11%  "
